# Accounts sheet gains a new data row (row 3): "test" in columns A-D,
# mirroring the header/User rows already present above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accounts")

$ws.Range("A3:D3").Value = "test"
